# Edit the resume document per the target diff.
$d = $word.ActiveDocument

# 1. Remove the " (current)" suffix after the Software Engineer Intern line.
$d.Content.Find.Execute(" (current)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null

# 2. Rewrite the bullet describing the unified analytics platform work.
$d.Content.Find.Execute(
    "Working on creating a unified, documented analytics platform for global Windows telemetry data",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Created a unified, documented data store for global Windows telemetry data specific to language usage",
    2) | Out-Null

# 3. Move the "_GoBack" bookmark from the "Developed internal tools" paragraph
#    to the end of the telemetry-data bullet we just rewrote.
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$target = $d.Content.Find
$target.Execute("specific to language usage", $true, $false, $false, $false,
                 $false, $true, 1, $false, "", 0) | Out-Null
$endRange = $d.Content
$endRange.Start = $target.Parent.End
$endRange.End = $target.Parent.End
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null
